# Apply regression-testing workbook update for the "2010-18" sheet:
#   - add a new run row (row 5, label "Baseline 2010-18 C369") below the
#     existing rows, duplicating the numbers from the "C267+" run (row 4)
#   - rename the row 4 run label from "C267+" to "C367+"
#   - leave the active selection on the new row's label cell (B5)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# 1. Populate the new row 5 first (its label becomes the newly-created
#    shared string) ...
$ws.Range("A5").Value = "CW3M"
$ws.Range("B5").Value = "Baseline 2010-18 C369"
$ws.Range("C5").Value = "2010-18"

$ws.Range("D5").Value = 677.97837322222222
$ws.Range("E5").Value = 2094.2995878888887
$ws.Range("F5").Value = 4.820043222222222
$ws.Range("G5").Value = 232.21855144444442
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 6.3389989999999994
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 586.02156566666656
$ws.Range("L5").Value = 96.631732222222212
$ws.Range("M5").Value = 1650.8734266666665
$ws.Range("N5").Value = 682.41798233333327
$ws.Range("O5").Value = 12820.605631666667
$ws.Range("P5").Value = 2216.7525497777779
$ws.Range("Q5").Value = 0.28915188888888882
$ws.Range("R5").Value = -0.0000024444444444444798

# 2. ... then rename the existing row 4 run label from "C267+" to "C367+".
$ws.Range("B4").Value = "Baseline 2010-18 C367+"

# 3. Leave the active selection on the new row's label cell, matching the
#    last user action.
$ws.Range("B5").Select()
